$p = $ppt.ActivePresentation

# --- Slide 2: "Activities this week" -----------------------------------
# The "Continuing to add functionality to Arduino communication classes"
# bullet used to be typed as two separate runs (" communication " and
# "classes"). Re-typing it collapses them into a single run with the same
# text, same formatting, same result on screen.
$s2 = $p.Slides.Item(2)
$contentShape = $s2.Shapes.Item(2)
$contentRange = $contentShape.TextFrame.TextRange
$arduinoPara = $contentRange.Paragraphs(3, 1)
$tail = $arduinoPara.Characters(43, 23)   # " communication classes"
$tail.Text = " communication classes"

# --- Slide 4: "Baseline Code Metrics" -----------------------------------
# The "# Source Files" metric was undercounted; correct it to 116 both in
# the on-slide summary table and in the linked trend chart beneath it.
$s4 = $p.Slides.Item(4)

$table = $s4.Shapes.Item(4).Table
$sourceFilesCell = $table.Cell(4, 2).Shape.TextFrame.TextRange
$sourceFilesCell.Characters(1, 3).Text = "116"

$chartShape = $s4.Shapes.Item(6)
$series = $chartShape.Chart.SeriesCollection().Item(1)
$series.Values = @(85, 85, 96, 97, 107, 112, 112, 112, 116)

# Line up the two trend charts at the same vertical position, and rename
# the second one now that it matches the first in height.
$chart1Shape = $s4.Shapes.Item(5)
$chart1Shape.Top = 309
$chartShape.Top = 309
$chartShape.Name = "Chart 8"
